$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 50.727272
$ws.Range("I38").Value = 50.727272
$ws.Range("K38").Value = 152.181816
$ws.Range("M38").Value = 219.818184
$ws.Range("H80").Value = 2393.4
$ws.Range("I80").Value = 2043
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 6129
$ws.Range("L80").Value = 8100
$ws.Range("M80").Value = -5131
$ws.Range("N80").Value = -10096
$ws.Range("H83").Value = 2393.4
$ws.Range("I83").Value = 2043
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 18387
$ws.Range("L83").Value = 24300
$ws.Range("M83").Value = -13395
$ws.Range("N83").Value = -34284
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()
$ws.Range("H111").Value = 13197.875
$ws.Range("I111").Value = 13197.875
$ws.Range("K111").Value = 39593.625
$ws.Range("M111").Value = -36526.625
$ws.Range("H113").Value = 76968.625
$ws.Range("I113").Value = 155792.86
$ws.Range("J113").Value = 15660.889
$ws.Range("K113").Value = 155792.86
$ws.Range("L113").Value = 15660.889
$ws.Range("M113").Value = -152538.86
$ws.Range("N113").Value = -22168.889
$ws.Range("H116").Value = 9696.458000000001
$ws.Range("I116").Value = 10864.667
$ws.Range("K116").Value = 10864.667
$ws.Range("M116").Value = -7422.666999999999
$ws.Range("H139").Value = 160000
$ws.Range("J139").Value = 180000
$ws.Range("L139").Value = 180000
$ws.Range("N139").Value = -190280
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 10000
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -9830
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -9730
$ws.Range("N20").ClearContents()
$ws.Range("H23").Value = 200000
$ws.Range("J23").Value = 200000
$ws.Range("L23").Value = 200000
$ws.Range("N23").Value = -200518
$ws.Range("H32").Value = 3860.0513
$ws.Range("I32").Value = 3195.5144
$ws.Range("K32").Value = 3195.5144
$ws.Range("M32").Value = -2908.5144
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21608
$ws.Range("H61").Value = 6028.778
$ws.Range("I61").Value = 5834.5
$ws.Range("K61").Value = 5834.5
$ws.Range("M61").Value = -5622.5
$ws.Range("H74").Value = 6777.1562
$ws.Range("J74").Value = 17673.363
$ws.Range("L74").Value = 17673.363
$ws.Range("N74").Value = -19421.363
$ws.Range("H77").Value = 6777.1562
$ws.Range("J77").Value = 17673.363
$ws.Range("L77").Value = 88366.815
$ws.Range("N77").Value = -97102.815
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws.Range("H130").Value = 101066.1
$ws.Range("J130").Value = 115499.375
$ws.Range("L130").Value = 115499.375
$ws.Range("N130").Value = -125539.375
$ws.Range("H132").Value = 4128.778
$ws.Range("I132").Value = 3498.3333
$ws.Range("J132").Value = 5389.6665
$ws.Range("K132").Value = 10494.9999
$ws.Range("L132").Value = 16168.9995
$ws.Range("M132").Value = -7964.999899999999
$ws.Range("N132").Value = -21228.9995
$ws.Range("H136").Value = 6028.778
$ws.Range("I136").Value = 5834.5
$ws.Range("K136").Value = 17503.5
$ws.Range("M136").Value = -14953.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 10243.286
$ws.Range("I11").Value = 23601.334
$ws.Range("J11").Value = 224.75
$ws.Range("K11").Value = 23601.334
$ws.Range("L11").Value = 224.75
$ws.Range("M11").Value = -23461.334
$ws.Range("N11").Value = -504.75
$ws.Range("H134").Value = 2183.139
$ws.Range("I134").Value = 1978.8064
$ws.Range("K134").Value = 5936.4192
$ws.Range("M134").Value = -3401.4192
$ws.Range("H139").Value = 89031
$ws.Range("J139").Value = 91288.75
$ws.Range("L139").Value = 91288.75
$ws.Range("N139").Value = -101568.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2999.5
$ws.Range("I17").Value = 2999.5
$ws.Range("K17").Value = 2999.5
$ws.Range("M17").Value = -2825.5
$ws.Range("H22").Value = 2420
$ws.Range("H31").Value = 31642.057
$ws.Range("I31").Value = 41260.32
$ws.Range("J31").Value = 7596.4
$ws.Range("K31").Value = 41260.32
$ws.Range("L31").Value = 7596.4
$ws.Range("M31").Value = -40965.32
$ws.Range("N31").Value = -8186.4
$ws.Range("H34").Value = 31642.057
$ws.Range("I34").Value = 41260.32
$ws.Range("J34").Value = 7596.4
$ws.Range("K34").Value = 41260.32
$ws.Range("L34").Value = 7596.4
$ws.Range("M34").Value = -41058.32
$ws.Range("N34").Value = -8000.4
$ws.Range("H44").Value = 42499.75
$ws.Range("J44").Value = 46666.332
$ws.Range("L44").Value = 46666.332
$ws.Range("N44").Value = -47550.332
$ws.Range("H55").Value = 12111.167
$ws.Range("I55").Value = 10533.4
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 10533.4
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -10218.4
$ws.Range("N55").Value = -20630
$ws.Range("H107").Value = 340.42856
$ws.Range("I107").Value = 340.42856
$ws.Range("K107").Value = 340.42856
$ws.Range("M107").Value = 1579.57144
$ws.Range("H137").Value = 93028.57000000001
$ws.Range("J137").Value = 93028.57000000001
$ws.Range("L137").Value = 93028.57000000001
$ws.Range("N137").Value = -103228.57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4072
$ws.Range("I22").Value = 900.5
$ws.Range("J22").Value = 4776.778
$ws.Range("K22").Value = 2701.5
$ws.Range("L22").Value = 14330.334
$ws.Range("M22").Value = -2532.5
$ws.Range("N22").Value = -14668.334
$ws.Range("H27").Value = 4072
$ws.Range("I27").Value = 900.5
$ws.Range("J27").Value = 4776.778
$ws.Range("K27").Value = 2701.5
$ws.Range("L27").Value = 14330.334
$ws.Range("M27").Value = -2599.5
$ws.Range("N27").Value = -14534.334
$ws.Range("H88").Value = 66667170
$ws.Range("J88").Value = 66667170
$ws.Range("L88").Value = 200001510
$ws.Range("N88").Value = -200002366
$ws.Range("H91").Value = 66667170
$ws.Range("J91").Value = 66667170
$ws.Range("L91").Value = 200001510
$ws.Range("N91").Value = -200004474
$ws.Range("H138").Value = 16675811
$ws.Range("I138").Value = 62511510
$ws.Range("K138").Value = 187534530
$ws.Range("M138").Value = -187529390
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5253.857
$ws.Range("J80").Value = 9356
$ws.Range("L80").Value = 9356
$ws.Range("N80").Value = -11352
$ws.Range("H83").Value = 5253.857
$ws.Range("J83").Value = 9356
$ws.Range("L83").Value = 46780
$ws.Range("N83").Value = -56764
$ws.Range("H113").Value = 1757.2727
$ws.Range("J113").Value = 1443
$ws.Range("L113").Value = 1443
$ws.Range("N113").Value = -5783
$ws.Range("H132").Value = 2991.4285
$ws.Range("I132").Value = 2823.5
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 8470.5
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -5940.5
$ws.Range("N132").Value = -17057
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1174.091
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 1174.091
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H140").Value = 60000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 60000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 60000
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -70360
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2376.3333
$ws.Range("I136").Value = 2574.8462
$ws.Range("K136").Value = 7724.5386
$ws.Range("M136").Value = -5174.5386
